$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 6.336
$ws.Range("A3").Value = -21.391
$ws.Range("B5").Value = 6.712999999999999
$ws.Range("D5").Value = -8.395
$ws.Range("D9").Value = -7.925
$ws.Range("D11").Value = -8.359
$ws.Range("A14").Value = -20.891
$ws.Range("A21").Value = -21.072
$ws.Range("D21").Value = -7.826000000000001
$ws.Range("A23").Value = -21.666
$ws.Range("A25").Value = -22.078
